$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "25/03/2023"
$ws.Range("D3").Value = 258
$ws.Range("D4").Value = 270
$ws.Range("D5").Value = 312
$ws.Range("D6").Value = 336
$ws.Range("D7").Value = 230
$ws.Range("D8").Value = 260
$ws.Range("D9").Value = 204
$ws.Range("D10").Value = 192
$ws.Range("D11").Value = 198
$ws.Range("D12").Value = 186
$ws.Range("D13").Value = 110
$ws.Range("D14").Value = 62
$ws.Range("D18").Value = 3
$ws.Range("D19").Value = 10
$ws.Range("D20").Value = 26
$ws.Range("D21").Value = 16
$ws.Range("D22").Value = 23
$ws.Range("D23").Value = 13
$ws.Range("D24").Value = 16
$ws.Range("D25").Value = 8
$ws.Range("D26").Value = 11
$ws.Range("D27").Value = 9
$ws.Range("D28").Value = 11
$ws.Range("D29").Value = 5
$ws.Range("D30").Value = 4
$ws.Range("D34").Value = 23
$ws.Range("D35").Value = 40
$ws.Range("D36").Value = 78
$ws.Range("D37").Value = 167
$ws.Range("D38").Value = 138
$ws.Range("D39").Value = 125
$ws.Range("D40").Value = 116
$ws.Range("D41").Value = 99
$ws.Range("D42").Value = 98
$ws.Range("D43").Value = 127
$ws.Range("D44").Value = 105
$ws.Range("D45").Value = 105
$ws.Range("D46").Value = 91
$ws.Range("D47").Value = 82
$ws.Range("D48").Value = 68
$ws.Range("D49").Value = 27
$ws.Range("D50").Value = 27
$ws.Range("D51").Value = 22
$ws.Range("D52").Value = 55
$ws.Range("D53").Value = 80
$ws.Range("D54").Value = 87
$ws.Range("D55").Value = 96
$ws.Range("D56").Value = 98
$ws.Range("D57").Value = 73
$ws.Range("D58").Value = 75
$ws.Range("D59").Value = 59
$ws.Range("D60").Value = 65
$ws.Range("D61").Value = 60
$ws.Range("D62").Value = 57
$ws.Range("D63").Value = 38
$ws.Range("D68").Value = 16
$ws.Range("D69").Value = 25
$ws.Range("D70").Value = 26
$ws.Range("D71").Value = 30
$ws.Range("D72").Value = 29
$ws.Range("D73").Value = 26
$ws.Range("D74").Value = 22
$ws.Range("D75").Value = 22
$ws.Range("D76").Value = 24
$ws.Range("D77").Value = 22
$ws.Range("D78").Value = 19
$ws.Range("D79").Value = 13
$ws.Range("D80").Value = 7
$ws.Range("D81").Value = 4
$ws.Range("D83").Value = 1
$ws.Range("D85").Value = 8
$ws.Range("D86").Value = 1
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 2
$ws.Range("D91").Value = 80
$ws.Range("D92").Value = 62
$ws.Range("D93").Value = 93
$ws.Range("D94").Value = 127
$ws.Range("D95").Value = 136
$ws.Range("D96").Value = 136
$ws.Range("D97").Value = 109
$ws.Range("D98").Value = 115
$ws.Range("D99").Value = 103
$ws.Range("D100").Value = 101
$ws.Range("D101").Value = 94
$ws.Range("D102").Value = 92
$ws.Range("D103").Value = 71
$ws.Range("D104").Value = 38
$ws.Range("D105").Value = 29
$ws.Range("D106").Value = 15
$ws.Range("D107").Value = 8
$ws.Range("D109").Value = 9
$ws.Range("D110").Value = 19
$ws.Range("D111").Value = 16
$ws.Range("D112").Value = 14
$ws.Range("D113").Value = 8
$ws.Range("D114").Value = 7
$ws.Range("D115").Value = 6
$ws.Range("D116").Value = 7
$ws.Range("D117").Value = 6
$ws.Range("D118").Value = 4
$ws.Range("D119").Value = 5
$ws.Range("D120").Value = 2
$ws.Range("D125").Value = 0
$ws.Range("D126").Value = 1
$ws.Range("D127").Value = 2
$ws.Range("D128").Value = 0
$ws.Range("D129").Value = 0
$ws.Range("D133").Value = 0
$ws.Range("D140").Value = 4
$ws.Range("D141").Value = 6
$ws.Range("D142").Value = 8
$ws.Range("D143").Value = 9
$ws.Range("D144").Value = 10
$ws.Range("D145").Value = 8
$ws.Range("D146").Value = 7
$ws.Range("D147").Value = 8
$ws.Range("D148").Value = 6
$ws.Range("D149").Value = 8
$ws.Range("D150").Value = 9
$ws.Range("D151").Value = 4
$ws.Range("D152").Value = 2
$ws.Range("D153").Value = 1
$ws.Range("D154").Value = 1
